$wb = $excel.ActiveWorkbook

# @@ -1468,25 +1468,25 @@
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2065.647
$ws.Range("J17").Value = 2186.8572
$ws.Range("L17").Value = 6560.571599999999
$ws.Range("N17").Value = -6896.571599999999

# @@ -6195,25 +6195,25 @@
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1204.0625
$ws.Range("J112").Value = 1212.7307
$ws.Range("L112").Value = 3638.1921
$ws.Range("N112").Value = -5854.1921

# @@ -6247,25 +6247,25 @@
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 9978.286
$ws.Range("I113").Value = 12402.7
$ws.Range("J113").Value = 7774.273
$ws.Range("K113").Value = 12402.7
$ws.Range("L113").Value = 7774.273
$ws.Range("M113").Value = -9148.700000000001
$ws.Range("N113").Value = -14282.273

# @@ -7502,25 +7502,25 @@
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1742.8682
$ws.Range("J138").Value = 1751.7595
$ws.Range("L138").Value = 5255.278499999999
$ws.Range("N138").Value = -15535.2785

# @@ -9292,25 +9292,25 @@
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6184.5674
$ws.Range("I32").Value = 3923.0908
$ws.Range("J32").Value = 16549.666
$ws.Range("K32").Value = 3923.0908
$ws.Range("L32").Value = 16549.666
$ws.Range("M32").Value = -3636.0908
$ws.Range("N32").Value = -17123.666

# @@ -11323,22 +11323,22 @@
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1959.3334
$ws.Range("I74").Value = 1641.3
$ws.Range("K74").Value = 1641.3
$ws.Range("M74").Value = -767.3

# @@ -11470,22 +11470,22 @@
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1959.3334
$ws.Range("I77").Value = 1641.3
$ws.Range("K77").Value = 8206.5
$ws.Range("M77").Value = -3838.5

# @@ -12441,22 +12441,22 @@
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 2403.4614
$ws.Range("I97").Value = 1393.1428
$ws.Range("K97").Value = 1393.1428
$ws.Range("M97").Value = -897.1428000000001

# @@ -13075,22 +13075,22 @@
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1349.75
$ws.Range("I110").Value = 1317.909
$ws.Range("K110").Value = 1317.909
$ws.Range("M110").Value = 727.0909999999999

# @@ -13666,22 +13666,22 @@
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 4323.609
$ws.Range("I122").Value = 4535.533
$ws.Range("K122").Value = 13606.599
$ws.Range("M122").Value = -11156.599

# @@ -13865,22 +13865,22 @@
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H126").Value = 5570555
$ws.Range("I126").Value = 5570555
$ws.Range("K126").Value = 16711665
$ws.Range("M126").Value = -16709195

# @@ -19646,19 +19646,22 @@
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H102").Value = 8888
$ws.Range("I102").Value = 8888
$ws.Range("K102").Value = 8888
$ws.Range("M102").Value = -5643

# @@ -20635,22 +20638,22 @@
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H122").Value = 70852
$ws.Range("J122").Value = 70852
$ws.Range("L122").Value = 70852
$ws.Range("N122").Value = -80652

# @@ -20972,22 +20975,22 @@
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H129").Value = 99973
$ws.Range("J129").Value = 99973
$ws.Range("L129").Value = 99973
$ws.Range("N129").Value = -109973

# @@ -21021,22 +21024,22 @@
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H130").Value = 101122.25
$ws.Range("J130").Value = 101122.25
$ws.Range("L130").Value = 101122.25
$ws.Range("N130").Value = -111162.25

# @@ -21070,19 +21073,22 @@
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H131").Value = 55000
$ws.Range("J131").Value = 55000
$ws.Range("L131").Value = 55000
$ws.Range("N131").Value = -65080

# @@ -23166,22 +23172,22 @@
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4019.1538
$ws.Range("I31").Value = 1925.3529
$ws.Range("K31").Value = 1925.3529
$ws.Range("M31").Value = -1630.3529

# @@ -23319,22 +23325,22 @@
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4019.1538
$ws.Range("I34").Value = 1925.3529
$ws.Range("K34").Value = 1925.3529
$ws.Range("M34").Value = -1723.3529

# @@ -24152,22 +24158,22 @@
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 47041.8
$ws.Range("I51").Value = 37403.168
$ws.Range("K51").Value = 37403.168
$ws.Range("M51").Value = -36667.168

# @@ -24654,22 +24660,22 @@
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 47041.8
$ws.Range("I61").Value = 37403.168
$ws.Range("K61").Value = 37403.168
$ws.Range("M61").Value = -37055.168

# @@ -26528,19 +26534,19 @@
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I99").Value = 2229.25
$ws.Range("K99").Value = 2229.25
$ws.Range("M99").Value = -731.25

# @@ -27652,25 +27658,25 @@
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 3006.15
$ws.Range("I122").Value = 2409.5386
$ws.Range("J122").Value = 4114.143
$ws.Range("K122").Value = 7228.6158
$ws.Range("L122").Value = 12342.429
$ws.Range("M122").Value = -4778.6158
$ws.Range("N122").Value = -17242.429

# @@ -27854,19 +27860,19 @@
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I126").Value = 2229.25
$ws.Range("K126").Value = 6687.75
$ws.Range("M126").Value = -4217.75

# @@ -28136,22 +28142,22 @@
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1874.7646
$ws.Range("I132").Value = 1441.8334
$ws.Range("K132").Value = 4325.5002
$ws.Range("M132").Value = -1795.5002

# @@ -28778,25 +28784,25 @@
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3469.7334
$ws.Range("I3").Value = 913.2727
$ws.Range("J3").Value = 10500
$ws.Range("K3").Value = 2739.8181
$ws.Range("L3").Value = 31500
$ws.Range("M3").Value = -2627.8181
$ws.Range("N3").Value = -31724

# @@ -28882,25 +28888,25 @@
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2141.6667
$ws.Range("J5").Value = 3042.6428
$ws.Range("L5").Value = 9127.928400000001
$ws.Range("N5").Value = -9351.928400000001

# @@ -31425,22 +31431,22 @@
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 6609.615
$ws.Range("J55").Value = 6609.615
$ws.Range("L55").Value = 19828.845
$ws.Range("N55").Value = -20182.845

# @@ -35423,25 +35429,25 @@
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 2141.6667
$ws.Range("J135").Value = 3042.6428
$ws.Range("L135").Value = 27383.7852
$ws.Range("N135").Value = -32453.7852

# @@ -38092,25 +38098,22 @@
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 50000
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 50000
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 50000
$ws.Range("M47").ClearContents()
$ws.Range("N47").Value = -51136

# @@ -38144,19 +38147,22 @@
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 6999.5
$ws.Range("J48").Value = 6999.5
$ws.Range("L48").Value = 6999.5
$ws.Range("N48").Value = -7969.5

# @@ -40775,22 +40781,22 @@
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1243.5555
$ws.Range("I102").Value = 1315.8334
$ws.Range("K102").Value = 1315.8334
$ws.Range("M102").Value = 306.1666

# @@ -41761,22 +41767,22 @@
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 19477.572
$ws.Range("I122").Value = 21140.5
$ws.Range("K122").Value = 63421.5
$ws.Range("M122").Value = -60971.5

# @@ -41960,22 +41966,22 @@
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 7777.3
$ws.Range("I126").Value = 3037.5
$ws.Range("K126").Value = 9112.5
$ws.Range("M126").Value = -6642.5

# @@ -42248,25 +42254,25 @@
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2589.5264
$ws.Range("I132").Value = 2583.0715
$ws.Range("J132").Value = 2607.6
$ws.Range("K132").Value = 7749.2145
$ws.Range("L132").Value = 7822.799999999999
$ws.Range("M132").Value = -5219.2145
$ws.Range("N132").Value = -12882.8

# @@ -43812,25 +43818,25 @@
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 918.1739
$ws.Range("I22").Value = 556.1818
$ws.Range("J22").Value = 1250
$ws.Range("K22").Value = 556.1818
$ws.Range("L22").Value = 1250
$ws.Range("M22").Value = -261.1818
$ws.Range("N22").Value = -1840

# @@ -44063,25 +44069,25 @@
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 918.1739
$ws.Range("I27").Value = 556.1818
$ws.Range("J27").Value = 1250
$ws.Range("K27").Value = 556.1818
$ws.Range("L27").Value = 1250
$ws.Range("M27").Value = -449.1818
$ws.Range("N27").Value = -1464

# @@ -44691,22 +44697,22 @@
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2417644.2
$ws.Range("I40").Value = 2047.8889
$ws.Range("K40").Value = 2047.8889
$ws.Range("M40").Value = -1911.8889

# @@ -45423,25 +45429,25 @@
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 2099.9707
$ws.Range("I55").Value = 1152.1305
$ws.Range("J55").Value = 4081.818
$ws.Range("K55").Value = 1152.1305
$ws.Range("L55").Value = 4081.818
$ws.Range("M55").Value = -979.1305
$ws.Range("N55").Value = -4427.818

# @@ -47463,22 +47469,19 @@
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()

# @@ -49166,22 +49169,22 @@
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3130.2778
$ws.Range("I132").Value = 1991.6666
$ws.Range("K132").Value = 5974.9998
$ws.Range("M132").Value = -3444.9998

# @@ -51921,22 +51924,22 @@
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 204664
$ws.Range("J46").Value = 204664
$ws.Range("L46").Value = 204664
$ws.Range("N46").Value = -205126

# @@ -52160,22 +52163,25 @@
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 30333.334
$ws.Range("I51").Value = 23000
$ws.Range("J51").Value = 45000
$ws.Range("K51").Value = 23000
$ws.Range("L51").Value = 45000
$ws.Range("M51").Value = -22490
$ws.Range("N51").Value = -46020

# @@ -56209,22 +56215,22 @@
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H134").Value = 204664
$ws.Range("J134").Value = 204664
$ws.Range("L134").Value = 613992
$ws.Range("N134").Value = -619062
